# Generate Report for Handback
# Populate the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime" /
# "Error Detail" columns for the 7466ecc0-7166-4dd4-bde1-7a8458ba0175 row on both the
# zh-cn and de-de sheets, now that a (stale) handback has come in for that item.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3f3aec364583f5850007f9a9de7f2bc6e9cb2852/e2e/7466ecc0-7166-4dd4-bde1-7a8458ba0175.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a743cc5b42388dd82499dc45c345a17b311495f/e2e/7466ecc0-7166-4dd4-bde1-7a8458ba0175.md."
$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a743cc5b42388dd82499dc45c345a17b311495f/e2e/7466ecc0-7166-4dd4-bde1-7a8458ba0175.md"
$mdName = "7466ecc0-7166-4dd4-bde1-7a8458ba0175.md"

function Update-HandbackRow($ws, $lang, $handbackDateTime) {

    $xlfName = "7466ecc0-7166-4dd4-bde1-7a8458ba0175.9b182e36bd8c74d2fe7708aa24426ed4059b81be.$lang.xlf"

    # I7: Latest Target File -> becomes a hyperlink to the handback md, like the other rows
    $i7 = $ws.Range("I7")
    $i7.Value = $mdName
    $ws.Hyperlinks.Add($i7, $latestUrl, "", "", $mdName) | Out-Null

    # Match the blue/underlined look used by the other hyperlink cells on this sheet
    # (set after Hyperlinks.Add, which otherwise stamps its own default style on top)
    $i7.Font.Underline = $true
    $i7.Font.Color = 15570276

    # J7: Latest Handback File
    $ws.Range("J7").Value = $xlfName

    # K7: Latest Handback DateTime
    $ws.Range("K7").Value = $handbackDateTime

    # P7: Error Detail
    $ws.Range("P7").Value = $errorDetail
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackRow $wsZhCn "zh-cn" "2016-08-24 04:52:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackRow $wsDeDe "de-de" "2016-08-24 04:53:02"
